$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = 35.23671050776777
$ws.Range("E2").Value = 36.15439588328721

$ws.Range("D3").Value = 34.95350722098058
$ws.Range("E3").Value = 35.10024149369163

$ws.Range("D4").Value = 34.93877611932735
$ws.Range("E4").Value = 34.90285597640286

$ws.Range("D5").Value = 35.32282634241693
$ws.Range("E5").Value = 35.35382334326526

$ws.Range("D6").Value = 34.95869030693223
$ws.Range("E6").Value = 35.11888985603599

$ws.Range("D7").Value = 35.03854759314339
$ws.Range("E7").Value = 34.99056870452858

$ws.Range("D8").Value = 34.59192476885546
$ws.Range("E8").Value = 34.5203374798529

$ws.Range("D9").Value = 35.39743386168312
$ws.Range("E9").Value = 35.39667167644611
